$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" (column D) values are written as text (matching the source
# workbook, which stores them as inline strings like "22.364.93" or
# "1.004" -- values that Excel would otherwise auto-convert to numbers).
# Forcing NumberFormat to text ("@") before the write keeps the literal
# string; resetting the Style to "Normal" afterwards drops the temporary
# text format so the cell keeps its original (default) style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '22.364.93'
$ws.Range('E2').Value = '  -0.15%  '
Set-TextValue $ws.Range('D3') '1.566.57'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.40%  '
Set-TextValue $ws.Range('D5') '1.004'
$ws.Range('E5').Value = '  +0.31%  '
Set-TextValue $ws.Range('D6') '290.42'
$ws.Range('E6').Value = '  +0.16%  '
Set-TextValue $ws.Range('D7') '0.3743'
$ws.Range('E7').Value = '  +0.78%  '
Set-TextValue $ws.Range('D8') '49.05'
$ws.Range('E8').Value = '  -0.24%  '
Set-TextValue $ws.Range('D9') '0.3374'
$ws.Range('E9').Value = '  -0.65%  '
Set-TextValue $ws.Range('D10') '0.07516'
$ws.Range('E10').Value = '  -1.84%  '
Set-TextValue $ws.Range('D11') '1.125'
$ws.Range('E11').Value = '  -3.69%  '
$ws.Range('E12').Value = '  +0.41%  '
Set-TextValue $ws.Range('D13') '20.78'
$ws.Range('E13').Value = '  -3.36%  '
Set-TextValue $ws.Range('D14') '5.908'
$ws.Range('E14').Value = '  -2.39%  '
Set-TextValue $ws.Range('D15') '6.872'
$ws.Range('E15').Value = '  -0.90%  '
Set-TextValue $ws.Range('D16') '1.562.65'
$ws.Range('E16').Value = '  -0.41%  '
Set-TextValue $ws.Range('D17') '0.00001113'
$ws.Range('E17').Value = '  -1.55%  '
Set-TextValue $ws.Range('D18') '89.48'
$ws.Range('E18').Value = '  -0.84%  '
Set-TextValue $ws.Range('D19') '0.06720'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  +0.33%  '
Set-TextValue $ws.Range('D21') '6.160'
$ws.Range('E21').Value = '  -1.49%  '
Set-TextValue $ws.Range('D22') '16.35'
$ws.Range('E22').Value = '  -1.20%  '
Set-TextValue $ws.Range('D23') '11.84'
$ws.Range('E23').Value = '  -1.61%  '
Set-TextValue $ws.Range('D24') '22.367.74'
$ws.Range('E24').Value = '  -0.09%  '
Set-TextValue $ws.Range('D25') '2.382'
$ws.Range('E25').Value = '  +0.97%  '
Set-TextValue $ws.Range('D26') '2.701'
$ws.Range('E26').Value = '  -4.48%  '
Set-TextValue $ws.Range('D27') '19.99'
$ws.Range('E27').Value = '  -0.83%  '
Set-TextValue $ws.Range('D28') '147.37'
$ws.Range('E28').Value = '  +1.26%  '
Set-TextValue $ws.Range('D29') '4.997'
$ws.Range('E29').Value = '  +0.25%  '
Set-TextValue $ws.Range('D30') '124.81'
$ws.Range('E30').Value = '  -0.64%  '
Set-TextValue $ws.Range('D31') '1.739.74'
$ws.Range('E31').Value = '  -0.21%  '
Set-TextValue $ws.Range('D32') '2.017'
$ws.Range('E32').Value = '  +0.42%  '
Set-TextValue $ws.Range('D33') '0.9790'
$ws.Range('E33').Value = '  -2.93%  '
Set-TextValue $ws.Range('D34') '5.951'
$ws.Range('E34').Value = '  -4.56%  '
Set-TextValue $ws.Range('D35') '9.911'
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D36') '0.08443'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D37') '1.409'
$ws.Range('E37').Value = '  +9.63%  '
Set-TextValue $ws.Range('D38') '0.02455'
$ws.Range('E38').Value = '  -3.37%  '
Set-TextValue $ws.Range('D39') '0.2267'
$ws.Range('E39').Value = '  -2.49%  '
Set-TextValue $ws.Range('D40') '0.06419'
$ws.Range('E40').Value = '  -0.13%  '
Set-TextValue $ws.Range('D41') '5.353'
$ws.Range('E41').Value = '  -3.47%  '
Set-TextValue $ws.Range('D42') '0.6226'
$ws.Range('E42').Value = '  -1.97%  '
Set-TextValue $ws.Range('D43') '10.95'
$ws.Range('E43').Value = '  -6.66%  '
$ws.Range('E44').Value = '  +0.29%  '
Set-TextValue $ws.Range('D45') '13.92'
$ws.Range('E45').Value = '  -1.78%  '
Set-TextValue $ws.Range('D46') '3.793'
$ws.Range('E46').Value = '  +0.86%  '
Set-TextValue $ws.Range('D47') '0.5830'
$ws.Range('E47').Value = '  -2.65%  '
Set-TextValue $ws.Range('D48') '2.047'
$ws.Range('E48').Value = '  -2.68%  '
Set-TextValue $ws.Range('D49') '1.248'
$ws.Range('E49').Value = '  -1.37%  '
Set-TextValue $ws.Range('D50') '123.96'
$ws.Range('E50').Value = '  -0.38%  '
Set-TextValue $ws.Range('D51') '0.07315'
$ws.Range('E51').Value = '  +0.50%  '
